$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
